$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '26.873.86'
$ws.Range('E2').Value = '  -1.32%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.807.79'
$ws.Range('E3').Value = '  -1.04%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.0000'
$ws.Range('E4').Value = '  -0.59%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '310.22'
$ws.Range('E5').Value = '  -1.02%  '
$ws.Range('E6').Value = '  -0.46%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.4455'
$ws.Range('E7').Value = '  +4.83%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.3675'
$ws.Range('E8').Value = '  -1.02%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.07358'
$ws.Range('E9').Value = '  +1.40%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.8560'
$ws.Range('E10').Value = '  -0.96%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '20.68'
$ws.Range('E11').Value = '  -1.93%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '1.803.47'
$ws.Range('E12').Value = '  -1.36%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '6.610'
$ws.Range('E13').Value = '  -1.84%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '92.42'
$ws.Range('E14').Value = '  +3.27%  '
$ws.Range('B15').Value = 'Polkadot'
$ws.Range('C15').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '5.313'
$ws.Range('E15').Value = '  -0.13%  '
$ws.Range('B16').Value = 'TRON'
$ws.Range('C16').Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '0.07073'
$ws.Range('E16').Value = '  -0.29%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '1.001'
$ws.Range('E17').Value = '  -0.63%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '0.000008737'
$ws.Range('E18').Value = '  -1.60%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '14.88'
$ws.Range('E20').Value = '  -1.48%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '26.890.08'
$ws.Range('E21').Value = '  -1.64%  '
$ws.Range('E22').Value = '  +0.46%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '10.87'
$ws.Range('E23').Value = '  -0.34%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '1.993'
$ws.Range('E24').Value = '  -0.02%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '151.73'
$ws.Range('E25').Value = '  -0.57%  '
$ws.Range('B26').Value = 'EthereumClassic'
$ws.Range('C26').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '18.57'
$ws.Range('E26').Value = '  +0.77%  '
$ws.Range('B27').Value = 'LidoDAOToken'
$ws.Range('C27').Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '2.192'
$ws.Range('E27').Value = '  +0.38%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '5.223'
$ws.Range('E28').Value = '  -0.36%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '116.62'
$ws.Range('E29').Value = '  +0.07%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '0.08842'
$ws.Range('E30').Value = '  -0.29%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '0.7533'
$ws.Range('E31').Value = '  -0.61%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '1.176'
$ws.Range('E32').Value = '  -1.84%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '2.933'
$ws.Range('E33').Value = '  +4.78%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '4.473'
$ws.Range('E34').Value = '  +0.27%  '
$ws.Range('E35').Value = '  -0.54%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '1.090'
$ws.Range('E36').Value = '  -2.43%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.01969'
$ws.Range('E37').Value = '  -0.44%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.05202'
$ws.Range('E38').Value = '  -1.28%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.5329'
$ws.Range('E39').Value = '  +5.19%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '2.865'
$ws.Range('E40').Value = '  -0.24%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '7.026'
$ws.Range('E41').Value = '  -4.27%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.1692'
$ws.Range('E42').Value = '  -0.46%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.5164'
$ws.Range('E43').Value = '  +8.57%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '8.461'
$ws.Range('E44').Value = '  -2.65%  '
$ws.Range('B45').Value = 'RenderToken'
$ws.Range('C45').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '1.975'
$ws.Range('E45').Value = '  +6.27%  '
$ws.Range('B46').Value = 'EnergySwap'
$ws.Range('C46').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '10.51'
$ws.Range('E46').Value = '  -1.67%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '105.51'
$ws.Range('E47').Value = '  -1.91%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '1.679'
$ws.Range('E48').Value = '  +0.35%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.9996'
$ws.Range('E49').Value = '  -0.50%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.06340'
$ws.Range('E50').Value = '  -0.94%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.9211'
$ws.Range('E51').Value = '  +0.57%  '
